$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (row 2), pushing the
# existing 150 data rows down by two (old row N -> new row N+2).
$ws.Rows.Item(2).Resize(2).Insert()

# Row 2 - new customer record
$ws.Cells.Item(2, 1).Value = "KH"
$ws.Cells.Item(2, 2).Value = 382
$ws.Cells.Item(2, 3).Value = "Thị Minh"
$ws.Cells.Item(2, 4).Value = "CẦN THƠ"
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = "0988903782"
$ws.Cells.Item(2, 6).Style = "Normal"
$ws.Cells.Item(2, 9).Value = 200000
$ws.Cells.Item(2, 10).Value = 0

# Row 3 - new customer record
$ws.Cells.Item(3, 1).Value = "KH"
$ws.Cells.Item(3, 2).Value = 377
$ws.Cells.Item(3, 3).Value = "Huỳnh Huyền Trân"
$ws.Cells.Item(3, 4).Value = "CẦN THƠ"
$ws.Cells.Item(3, 9).Value = 6000000
$ws.Cells.Item(3, 10).Value = 0
